$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to repull/recalculated data
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -1
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = 1
$ws.Range("F16").Value = -1
$ws.Range("F20").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("F29").Value = 4
$ws.Range("F31").Value = -1
$ws.Range("F35").Value = 2
$ws.Range("F39").Value = 1
$ws.Range("F41").Value = 6
$ws.Range("F43").Value = -6
$ws.Range("F49").Value = -4
$ws.Range("F55").Value = -3
$ws.Range("F58").Value = 1
$ws.Range("F61").Value = 3
